$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.395.49'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '1.572.81'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '290.84'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3765'
$ws.Range('E7').Value = '  +2.41%  '
$ws.Range('E8').Value = '  +1.22%  '
$ws.Range('E9').Value = '  +1.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07647'
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('E12').Value = '  +0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.25'
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.031'
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.949'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value = '1.574.07'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001135'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '90.17'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06747'
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.86'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('D24').Value = '22.397.15'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.395'
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.691'
$ws.Range('E26').Value = '  -10.06%  '
$ws.Range('E27').Value = '  +1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '147.34'
$ws.Range('E28').Value = '  +1.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.028'
$ws.Range('E29').Value = '  +0.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.24'
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('D31').Value = '1.748.49'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.173'
$ws.Range('E32').Value = '  -1.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.007'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9907'
$ws.Range('E34').Value = '  -4.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.985'
$ws.Range('E35').Value = '  -3.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08530'
$ws.Range('E36').Value = '  +1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02554'
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06574'
$ws.Range('E39').Value = '  +0.53%  '
$ws.Range('E40').Value = '  +6.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.438'
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.52'
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6417'
$ws.Range('E43').Value = '  +0.27%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.15'
$ws.Range('E44').Value = '  -1.78%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6001'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.783'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.308'
$ws.Range('E48').Value = '  +6.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.094'
$ws.Range('E49').Value = '  -2.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '125.67'
$ws.Range('E50').Value = '  +2.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07330'
$ws.Range('E51').Value = '  +0.45%  '
